$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.131.28"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "2.471.74"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'519.50"
$ws.Range("E5").Value = "  -3.10%  "
$ws.Range("D6").Value = "'132.73"
$ws.Range("E6").Value = "  -3.78%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.558"
$ws.Range("D9").Value = "'0.0992"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "2.909.88"
$ws.Range("D14").Value = "58.058.32"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "'22.04"
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("D17").Value = "2.471.85"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "'10.85"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "'320.61"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").Value = "'64.34"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").Value = "'0.408"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'0.161"
$ws.Range("E26").Value = "  -3.53%  "
$ws.Range("D27").Value = "'7.39"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").Value = "'6.37"
$ws.Range("E29").Value = "  -5.27%  "
$ws.Range("D30").Value = "'1.70"
$ws.Range("E30").Value = "  -4.84%  "
$ws.Range("D31").Value = "'166.59"
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "'18.13"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("E36").Value = "  -10.28%  "
$ws.Range("D37").Value = "'3.99"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("D39").Value = "'0.793"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("D40").Value = "'276.87"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").Value = "'3.46"
$ws.Range("E41").Value = "  -4.87%  "
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("D43").Value = "'0.595"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").Value = "'126.04"
$ws.Range("E44").Value = "  -4.82%  "
$ws.Range("D45").Value = "'0.0908"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "'0.0493"
$ws.Range("E46").Value = "  -3.48%  "
$ws.Range("D47").Value = "'0.0214"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").Value = "'17.10"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").Value = "1.735.87"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "'4.67"
$ws.Range("E51").Value = "  -1.75%  "
